# Change figures and documentation
#   - "removegroup"  -> "deletegroup"   (lower-case label, e.g. shape "TextBox 18")
#   - "removeGroup"  -> "deleteGroup"   (camel-case label, e.g. shape "TextBox 22")
#
# Each occurrence is replaced in place via TextRange.Characters(start,length),
# which addresses an exact character span without disturbing neighbouring
# runs or any run-level formatting (color, size, the spell-check "err"
# flag, ...). This keeps e.g. "removegroup" + " " as two separate runs,
# and "removeGroup" + "()" as two separate runs, exactly as authored.

$p = $ppt.ActivePresentation

function Replace-Literal {
    param(
        [string]$oldWord,
        [string]$newWord
    )

    if ($oldWord.Length -eq 0) { return }

    for ($si = 1; $si -le $p.Slides.Count; $si++) {
        $slide = $p.Slides.Item($si)
        for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
            $shape = $slide.Shapes.Item($shi)
            if (-not $shape.HasTextFrame) { continue }
            if (-not $shape.TextFrame.HasText) { continue }

            $tr = $shape.TextFrame.TextRange
            $text = $tr.Text

            $searchFrom = 0
            while ($searchFrom -le $text.Length) {
                $idx = $text.IndexOf($oldWord, $searchFrom)
                if ($idx -lt 0) { break }

                # PowerPoint TextRange character positions are 1-based.
                $span = $tr.Characters($idx + 1, $oldWord.Length)
                $span.Text = $newWord

                # Refresh the cached text / continue scanning after the
                # replacement (lengths may differ between old/new words).
                $text = $tr.Text
                $searchFrom = $idx + $newWord.Length
            }
        }
    }
}

Replace-Literal "removegroup" "deletegroup"
Replace-Literal "removeGroup" "deleteGroup"
